# Project Backlog.xlsx - fill in "Initial Estimate" / "Adjust Factor" data
# for every backlog item on Sheet1. The "Adjust Estimate" column (E) already
# holds the formula =Cn*(1+Dn) for each row, plus SUM() totals in row 28, so
# simply populating C/D lets the workbook's own formulas recompute E and the
# totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row -> Initial Estimate (column C). Adjust Factor (column D) is 30% (0.3)
# for every row, including the totals row.
$estimates = [ordered]@{
    2  = 6
    3  = 6
    4  = 6
    6  = 6
    8  = 6
    9  = 6
    11 = 4
    12 = 4
    13 = 8
    15 = 3
    16 = 4
    17 = 4
    18 = 8
    20 = 6
    21 = 4
    23 = 4
    24 = 8
    26 = 4
}

foreach ($row in $estimates.Keys) {
    $ws.Cells.Item($row, 3).Value = $estimates[$row]   # column C
    $ws.Cells.Item($row, 4).Value = 0.3                # column D
}

# Totals row: C28/E28 are already SUM() formulas, but it also carries an
# Adjust Factor value of 0.3.
$ws.Cells.Item(28, 4).Value = 0.3

# Recalculate so every dependent formula (E column + totals) picks up the
# new inputs.
$excel.Calculate() | Out-Null

# Restore the view/selection state: scrolled so row 22 is at the top, with
# the active cell/selection on I14.
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("I14").Select() | Out-Null
